$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) from Excel auto-converting numeric-looking text
# (e.g. "1.005", "28.692.86") into actual numbers. Apply a Text number format
# to the whole price column before writing values, then clear the format again
# so the cell style matches the original workbook (no explicit style index).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Price column (D) updates ---
$ws.Range("D2").Value = '28.692.86'
$ws.Range("D3").Value = '1.869.74'
$ws.Range("D4").Value = '1.006'
$ws.Range("D5").Value = '327.29'
$ws.Range("D6").Value = '1.005'
$ws.Range("D7").Value = '0.4637'
$ws.Range("D8").Value = '0.3917'
$ws.Range("D9").Value = '0.07924'
$ws.Range("D10").Value = '0.9702'
$ws.Range("D11").Value = '22.30'
$ws.Range("D12").Value = '1.902.24'
$ws.Range("D13").Value = '5.735'
$ws.Range("D14").Value = '6.940'
$ws.Range("D15").Value = '0.06963'
$ws.Range("D16").Value = '88.27'
$ws.Range("D18").Value = '0.00001007'
$ws.Range("D20").Value = '1.006'
$ws.Range("D21").Value = '28.702.98'
$ws.Range("D22").Value = '5.329'
$ws.Range("D24").Value = '2.125'
$ws.Range("D25").Value = '2.143.89'
$ws.Range("D26").Value = '153.55'
$ws.Range("D27").Value = '19.34'
$ws.Range("D28").Value = '5.712'
$ws.Range("D29").Value = '2.004'
$ws.Range("D30").Value = '119.62'
$ws.Range("D31").Value = '0.09372'
$ws.Range("D32").Value = '0.9322'
$ws.Range("D33").Value = '5.330'
$ws.Range("D34").Value = '1.346'
$ws.Range("D35").Value = '3.359'
$ws.Range("D36").Value = '0.05843'
$ws.Range("D37").Value = '0.02130'
$ws.Range("D38").Value = '1.151'
$ws.Range("D39").Value = '7.912'
$ws.Range("D41").Value = '9.956'
$ws.Range("D43").Value = '0.07246'
$ws.Range("D44").Value = '11.82'
$ws.Range("D45").Value = '0.5324'
$ws.Range("D47").Value = '1.139'
$ws.Range("D48").Value = '1.849'
$ws.Range("D50").Value = '2.350'
$ws.Range("D51").Value = '1.005'

# Restore the default (no explicit number format) style on the price column
$priceRange.ClearFormats()

# --- Other column (B, C, E) updates ---
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("E10").Value = '  +0.56%  '
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("E12").Value = '  +5.15%  '
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("E25").Value = '  +3.57%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("E30").Value = '  +2.16%  '
$ws.Range("E32").Value = '  -1.44%  '
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("E35").Value = '  -2.42%  '
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  +3.41%  '
$ws.Range("E40").Value = '  +0.53%  '
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E43").Value = '  +2.88%  '
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  -5.36%  '
$ws.Range("E47").Value = '  -7.90%  '
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("E51").Value = '  +0.49%  '

Write-Host "Applied crypto price/volume update."
